# "Generate Report for Handback"
# Marks the two tracked files (53478499-... and 839848ef-...) as handed back:
# updates the Status text everywhere it appears, fills in the "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" columns on the
# zh-cn and de-de detail sheets (with a hyperlink on the target-file cell,
# like the existing Source File Name column), and widens a few columns that
# now hold longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

$newStatus = "Handed back: in sync with en-US"

$mdUrl53478499 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3667e744e7d99241111e92340227cf2336251ff6/e2e/53478499-934e-480e-bd4b-871f57b73d9c.md"
$mdUrl839848ef = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3667e744e7d99241111e92340227cf2336251ff6/e2e/839848ef-74c6-441f-b782-c0f541aa1a71.md"

$md53478499 = "53478499-934e-480e-bd4b-871f57b73d9c.md"
$md839848ef = "839848ef-74c6-441f-b782-c0f541aa1a71.md"

# ---------------------------------------------------------------------------
# 1. Status flips from "Ready for handoff" to "Handed back: in sync with
#    en-US" everywhere it is shown: the Overview per-locale columns and the
#    Status column on each locale detail sheet.
# ---------------------------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn detail sheet: fill in Latest Target File / Latest Handback File /
#    Latest Handback DateTime for both rows.
# ---------------------------------------------------------------------------
$zhcn.Range("I2").Value = $md53478499
$zhcn.Range("J2").Value = "53478499-934e-480e-bd4b-871f57b73d9c.e2160e0be45c77815f671bc1b7d30101bbb330bc.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-05 13:16:33"

$zhcn.Range("I3").Value = $md839848ef
$zhcn.Range("J3").Value = "839848ef-74c6-441f-b782-c0f541aa1a71.7432e153b21e08bcc6d0056b3d9978990a763004.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-05 13:16:33"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl53478499, [Type]::Missing, [Type]::Missing, $md53478499) | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdUrl839848ef, [Type]::Missing, [Type]::Missing, $md839848ef) | Out-Null

# ---------------------------------------------------------------------------
# 3. de-de detail sheet: same shape, different xlf/datetime values.
# ---------------------------------------------------------------------------
$dede.Range("I2").Value = $md53478499
$dede.Range("J2").Value = "53478499-934e-480e-bd4b-871f57b73d9c.e2160e0be45c77815f671bc1b7d30101bbb330bc.de-de.xlf"
$dede.Range("K2").Value = "2016-09-05 13:16:41"

$dede.Range("I3").Value = $md839848ef
$dede.Range("J3").Value = "839848ef-74c6-441f-b782-c0f541aa1a71.7432e153b21e08bcc6d0056b3d9978990a763004.de-de.xlf"
$dede.Range("K3").Value = "2016-09-05 13:16:41"

$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl53478499, [Type]::Missing, [Type]::Missing, $md53478499) | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), $mdUrl839848ef, [Type]::Missing, [Type]::Missing, $md839848ef) | Out-Null

# ---------------------------------------------------------------------------
# 4. Column widths: the Status / Target-File / Handback-File columns now
#    carry longer strings, so they are widened to fit. Excel's ColumnWidth
#    is stored in whole-pixel steps, so we feed values that land on the
#    nearest achievable step to the intended width.
# ---------------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.16666667
$overview.Columns.Item(6).ColumnWidth = 29.16666667

$zhcn.Columns.Item(3).ColumnWidth = 29.16666667
$zhcn.Columns.Item(9).ColumnWidth = 39.16666667
$zhcn.Columns.Item(10).ColumnWidth = 39.16666667

$dede.Columns.Item(3).ColumnWidth = 29.16666667
$dede.Columns.Item(9).ColumnWidth = 39.16666667
$dede.Columns.Item(10).ColumnWidth = 39.16666667

Write-Output "Handback report generated"
